$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update price (D) and volume-change (E) columns for rows with refreshed crypto data ---
$ws.Cells.Item(2, 4).Value = "'27.379.64"
$ws.Cells.Item(2, 5).Value = "  -2.25%  "
$ws.Cells.Item(3, 4).Value = "'1.710.58"
$ws.Cells.Item(3, 5).Value = "  -1.74%  "
$ws.Cells.Item(4, 4).Value = "'1.005"
$ws.Cells.Item(4, 5).Value = "  +0.29%  "
$ws.Cells.Item(5, 4).Value = "'224.09"
$ws.Cells.Item(5, 5).Value = "  -1.95%  "
$ws.Cells.Item(6, 4).Value = "'0.5325"
$ws.Cells.Item(6, 5).Value = "  -2.04%  "
$ws.Cells.Item(7, 4).Value = "'1.004"
$ws.Cells.Item(7, 5).Value = "  +0.20%  "
$ws.Cells.Item(8, 4).Value = "'0.2661"
$ws.Cells.Item(8, 5).Value = "  -3.83%  "
$ws.Cells.Item(9, 4).Value = "'0.06595"
$ws.Cells.Item(9, 5).Value = "  -1.87%  "
$ws.Cells.Item(10, 4).Value = "'20.89"
$ws.Cells.Item(10, 5).Value = "  -3.59%  "
$ws.Cells.Item(11, 4).Value = "'0.07642"
$ws.Cells.Item(11, 5).Value = "  -1.87%  "
$ws.Cells.Item(12, 4).Value = "'4.574"
$ws.Cells.Item(12, 5).Value = "  -2.59%  "
$ws.Cells.Item(13, 4).Value = "'1.716.24"
$ws.Cells.Item(13, 5).Value = "  -1.62%  "
$ws.Cells.Item(14, 4).Value = "'1.945.15"
$ws.Cells.Item(14, 5).Value = "  -1.72%  "
$ws.Cells.Item(15, 4).Value = "'0.5729"
$ws.Cells.Item(15, 5).Value = "  -4.01%  "
$ws.Cells.Item(16, 4).Value = "'0.0₅8171"
$ws.Cells.Item(16, 5).Value = "  -2.41%  "
$ws.Cells.Item(17, 4).Value = "'67.86"
$ws.Cells.Item(17, 5).Value = "  -1.61%  "
$ws.Cells.Item(18, 4).Value = "'27.361.19"
$ws.Cells.Item(18, 5).Value = "  -2.24%  "
$ws.Cells.Item(19, 4).Value = "'216.12"
$ws.Cells.Item(19, 5).Value = "  -3.78%  "
$ws.Cells.Item(20, 5).Value = "  +0.17%  "
$ws.Cells.Item(21, 4).Value = "'4.673"
$ws.Cells.Item(21, 5).Value = "  -3.37%  "
$ws.Cells.Item(22, 4).Value = "'10.43"
$ws.Cells.Item(22, 5).Value = "  -4.17%  "
$ws.Cells.Item(23, 4).Value = "'5.974"
$ws.Cells.Item(23, 5).Value = "  -4.04%  "
$ws.Cells.Item(24, 5).Value = "  +0.04%  "
$ws.Cells.Item(25, 4).Value = "'1.766"
$ws.Cells.Item(25, 5).Value = "  +7.05%  "
$ws.Cells.Item(26, 4).Value = "'141.64"
$ws.Cells.Item(26, 5).Value = "  -3.15%  "
$ws.Cells.Item(27, 4).Value = "'0.1218"
$ws.Cells.Item(27, 5).Value = "  -2.32%  "
$ws.Cells.Item(28, 4).Value = "'7.275"
$ws.Cells.Item(28, 5).Value = "  -2.51%  "
$ws.Cells.Item(29, 4).Value = "'16.32"
$ws.Cells.Item(29, 5).Value = "  -4.90%  "
$ws.Cells.Item(30, 4).Value = "'0.05417"
$ws.Cells.Item(30, 5).Value = "  -4.50%  "
$ws.Cells.Item(31, 4).Value = "'1.294"
$ws.Cells.Item(31, 5).Value = "  -1.88%  "
$ws.Cells.Item(32, 4).Value = "'3.507"
$ws.Cells.Item(32, 5).Value = "  -5.56%  "
$ws.Cells.Item(33, 4).Value = "'3.430"
$ws.Cells.Item(33, 5).Value = "  -2.66%  "
$ws.Cells.Item(34, 4).Value = "'1.645"
$ws.Cells.Item(34, 5).Value = "  -1.76%  "
$ws.Cells.Item(35, 4).Value = "'2.881"
$ws.Cells.Item(35, 5).Value = "  +0.72%  "
$ws.Cells.Item(38, 4).Value = "'0.5865"
$ws.Cells.Item(38, 5).Value = "  -1.53%  "
$ws.Cells.Item(39, 4).Value = "'0.01633"
$ws.Cells.Item(39, 5).Value = "  -2.41%  "
$ws.Cells.Item(40, 4).Value = "'5.873"
$ws.Cells.Item(40, 5).Value = "  -1.53%  "
$ws.Cells.Item(41, 4).Value = "'1.047.48"
$ws.Cells.Item(41, 5).Value = "  -0.19%  "
$ws.Cells.Item(42, 4).Value = "'1.004"
$ws.Cells.Item(42, 5).Value = "  +0.15%  "
$ws.Cells.Item(43, 4).Value = "'0.8411"
$ws.Cells.Item(43, 5).Value = "  -0.87%  "
$ws.Cells.Item(44, 4).Value = "'100.82"
$ws.Cells.Item(44, 5).Value = "  -1.28%  "
$ws.Cells.Item(45, 4).Value = "'1.851.07"
$ws.Cells.Item(45, 5).Value = "  -1.78%  "
$ws.Cells.Item(46, 5).Value = "  -0.39%  "
$ws.Cells.Item(47, 4).Value = "'57.98"
$ws.Cells.Item(47, 5).Value = "  -3.28%  "
$ws.Cells.Item(48, 4).Value = "'0.4508"
$ws.Cells.Item(48, 5).Value = "  +1.74%  "
$ws.Cells.Item(49, 4).Value = "'1.004"
$ws.Cells.Item(49, 5).Value = "  +0.65%  "
$ws.Cells.Item(50, 4).Value = "'8.101"
$ws.Cells.Item(50, 5).Value = "  -2.10%  "
$ws.Cells.Item(51, 4).Value = "'0.05246"
$ws.Cells.Item(51, 5).Value = "  -1.47%  "

# --- Row 36/37: HuobiToken and ARBITRUM swap places (with refreshed price/volume data) ---
$ws.Cells.Item(36, 2).Value = "ARBITRUM"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(36, 4).Value = "'0.9494"
$ws.Cells.Item(36, 5).Value = "  -3.29%  "

$ws.Cells.Item(37, 2).Value = "HuobiToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(37, 4).Value = "'2.419"
$ws.Cells.Item(37, 5).Value = "  -1.34%  "
